$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.332.13'
$ws.Range('E2').Value = '  -0.91%  '
$ws.Range('D3').Value = '2.772.85'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '352.49'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '107.95'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.549'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.19%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -1.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.58'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.83%  '
$ws.Range('E11').Value = '  +2.89%  '
$ws.Range('E12').Value = '  -1.76%  '
$ws.Range('E13').Value = '  +3.41%  '
$ws.Range('E14').Value = '  -0.55%  '
$ws.Range('D15').Value = '3.210.10'
$ws.Range('E15').Value = '  +0.21%  '
$ws.Range('D16').Value = '2.774.07'
$ws.Range('E16').Value = '  +0.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.920'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.04%  '
$ws.Range('D18').Value = '51.332.99'
$ws.Range('E18').Value = '  -0.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.59'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.48%  '
$ws.Range('E20').Value = '  -1.77%  '
$ws.Range('E21').Value = '  +1.02%  '
$ws.Range('E22').Value = '  -1.38%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.82'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '265.31'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.18%  '
$ws.Range('E25').Value = '  -0.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '26.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.07%  '
$ws.Range('E28').Value = '  +11.89%  '
$ws.Range('E29').Value = '  +0.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.41'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.11%  '
$ws.Range('E31').Value = '  -0.84%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.16'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +8.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '51.93'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.56%  '
$ws.Range('E34').Value = '  -2.29%  '
$ws.Range('E35').Value = '  +5.09%  '
$ws.Range('E36').Value = '  -2.08%  '
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.37'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.01%  '
$ws.Range('E39').Value = '  -2.32%  '
$ws.Range('E40').Value = '  -1.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.53'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.01%  '
$ws.Range('E42').Value = '  -0.61%  '
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '119.98'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.91%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.02'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.21%  '
$ws.Range('E45').Value = '  -2.33%  '
$ws.Range('D46').Value = '2.102.76'
$ws.Range('E46').Value = '  +1.96%  '
$ws.Range('E48').Value = '  +4.74%  '
$ws.Range('E49').Value = '  -5.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.902'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.72%  '
$ws.Range('E51').Value = '  +8.80%  '
